$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel;
# force Text format, assign, then restore the default "Normal" style so the
# cell keeps its original (unstyled) appearance while the stored value stays text.
$ws.Range('D2').Value = '63.698.17'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '3.117.67'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.113.25'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.516'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.147'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.31'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000250'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('D15').Value = '3.634.38'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').Value = '63.780.80'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '3.124.13'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '479.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.76%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -2.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.18'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.02'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.16%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  -8.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.71%  '
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.96'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0749'
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '435.12'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.77%  '
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.118'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range('D44').Value = '2.848.37'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.257'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.112'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.51'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.21%  '
